# agregar tiempo a la ventana checklist
# Applies the target edit to Hoja1:
#  - removes the duplicate "prueba"/test data row (old row 3)
#  - renames the last header to "Resultado Aprobador" (drops the trailing " 1")
#  - updates the Promocion/Numero Propuesta/Resultado columns on the remaining data row
#  - adds a new "Resultado Aprobador" value column (W) mirroring the emission result
#  - left-aligns the Promocion cell and text-formats the result columns
#  - widens the columns that now hold longer header/content text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the extra "prueba" data row (old row 3) entirely.
$ws.Rows.Item(3).Delete()

# 2) Header row: "Resultado Aprobador 1" -> "Resultado Aprobador"
$ws.Range("W1").Value = "Resultado Aprobador"

# 3) Data row (row 2) content updates
$ws.Range("F2").Value = 0
$ws.Range("U2").Value = "4900025"
$ws.Range("V2").Value = "El Documento ha sido derivado satisfactoriamente "
$ws.Range("W2").Value = "El Documento ha sido derivado satisfactoriamente "

# 4) Style adjustments
#    Promocion (F2) becomes a plain left-aligned cell
$ws.Range("F2").HorizontalAlignment = -4131
#    MontoAmortizar / Tasa Preferencial lose their old "text" number format
$ws.Range("I2").Style = "Normal"
$ws.Range("O2").Style = "Normal"
#    Numero Propuesta / Resultado* columns gain the "text" number format
$ws.Range("U2").NumberFormat = "@"
$ws.Range("V2").NumberFormat = "@"
$ws.Range("W2").NumberFormat = "@"

# 5) Column widths for the now-visible/resized columns
$ws.Columns.Item(2).ColumnWidth = 21.43
$ws.Columns.Item(3).ColumnWidth = 15.86
$ws.Columns.Item(4).ColumnWidth = 16.57
$ws.Columns.Item(5).ColumnWidth = 17.71
$ws.Columns.Item(8).ColumnWidth = 14.43
$ws.Columns.Item(9).ColumnWidth = 15.57
$ws.Columns.Item(11).ColumnWidth = 22.86
$ws.Columns.Item(15).ColumnWidth = 16.29
$ws.Columns.Item(17).ColumnWidth = 17.43
$ws.Columns.Item(18).ColumnWidth = 25.57
$ws.Columns.Item(19).ColumnWidth = 12
$ws.Columns.Item(21).ColumnWidth = 17.86
$ws.Columns.Item(22).ColumnWidth = 27.57
$ws.Columns.Item(23).ColumnWidth = 19.86

Write-Output "edit complete"
